$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 188: Nb nouveaux cas positifs (C188) corrected from 5 to 6
$ws.Range("C188").Value = 6

# Row 195: Nb nouveaux cas positifs (C195) corrected from 6 to 9
$ws.Range("C195").Value = 9

# Row 196: corrections
$ws.Range("C196").Value = 3
$ws.Range("D196").Value = 0
$ws.Range("I196").Value = 0

# Row 197: newly filled-in data (previously blank placeholder row)
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 0
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 4
$ws.Range("I197").Value = 0
$ws.Range("L197").Value = "0"
$ws.Range("M197").Value = "0"

# Scroll / selection state update to match the saved view
$win = $excel.ActiveWindow
$win.ScrollRow = 183
$win.ScrollColumn = 1
$ws.Range("C199").Select()
